$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 860.2143
$ws.Range("I41").Value = 606.8570999999999
$ws.Range("J41").Value = 1113.5714
$ws.Range("K41").Value = 606.8570999999999
$ws.Range("L41").Value = 1113.5714
$ws.Range("M41").Value = -166.8570999999999
$ws.Range("N41").Value = -1993.5714
$ws.Range("H53").Value = 5119
$ws.Range("I53").Value = 403.25
$ws.Range("J53").Value = 10777.9
$ws.Range("K53").Value = 403.25
$ws.Range("L53").Value = 10777.9
$ws.Range("M53").Value = 233.75
$ws.Range("N53").Value = -12051.9
$ws.Range("H88").Value = 2985.606
$ws.Range("I88").Value = 997.5
$ws.Range("J88").Value = 3113.8708
$ws.Range("K88").Value = 997.5
$ws.Range("L88").Value = 3113.8708
$ws.Range("M88").Value = -591.5
$ws.Range("N88").Value = -3925.8708
$ws.Range("H91").Value = 2985.606
$ws.Range("I91").Value = 997.5
$ws.Range("J91").Value = 3113.8708
$ws.Range("K91").Value = 997.5
$ws.Range("L91").Value = 3113.8708
$ws.Range("M91").Value = 406.5
$ws.Range("N91").Value = -5921.870800000001
$ws.Range("H106").Value = 83335830
$ws.Range("I106").Value = 100002296
$ws.Range("K106").Value = 100002296
$ws.Range("M106").Value = -100001665
$ws.Range("H112").Value = 7675.294
$ws.Range("J112").Value = 9120
$ws.Range("L112").Value = 27360
$ws.Range("N112").Value = -29576
$ws.Range("H131").Value = 5469.6
$ws.Range("I131").Value = 1939.2222
$ws.Range("J131").Value = 8358.091
$ws.Range("K131").Value = 5817.6666
$ws.Range("L131").Value = 25074.273
$ws.Range("M131").Value = -777.6665999999996
$ws.Range("N131").Value = -35154.273

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 6683
$ws.Range("I63").Value = 2378.5
$ws.Range("K63").Value = 2378.5
$ws.Range("M63").Value = -1692.5
$ws.Range("H66").Value = 6683
$ws.Range("I66").Value = 2378.5
$ws.Range("K66").Value = 11892.5
$ws.Range("M66").Value = -8460.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4168754
$ws.Range("I105").Value = 5210542.5
$ws.Range("K105").Value = 5210542.5
$ws.Range("M105").Value = -5208795.5
$ws.Range("H107").Value = 3573762.8
$ws.Range("I107").Value = 4203451.5
$ws.Range("K107").Value = 4203451.5
$ws.Range("M107").Value = -4201531.5
$ws.Range("H124").Value = 44409
$ws.Range("I124").Value = 44409
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 44409
$ws.Range("L124").Value = 0
$ws.Range("M124").ClearContents()
$ws.Range("N124").Value = -39499

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 1140.5
$ws.Range("I11").Value = 50
$ws.Range("J11").Value = 5502.5
$ws.Range("K11").Value = 50
$ws.Range("L11").Value = 5502.5
$ws.Range("M11").Value = 90
$ws.Range("N11").Value = -5782.5
$ws.Range("H58").Value = 6229.6113
$ws.Range("J58").Value = 3415.5386
$ws.Range("L58").Value = 3415.5386
$ws.Range("N58").Value = -3821.5386
$ws.Range("H132").Value = 83718.914
$ws.Range("I132").Value = 53822.527
$ws.Range("J132").Value = 225726.75
$ws.Range("K132").Value = 161467.581
$ws.Range("L132").Value = 677180.25
$ws.Range("M132").Value = -158937.581
$ws.Range("N132").Value = -682240.25
$ws.Range("H136").Value = 6229.6113
$ws.Range("J136").Value = 3415.5386
$ws.Range("L136").Value = 10246.6158
$ws.Range("N136").Value = -15346.6158
$ws.Range("H141").Value = 83177.46000000001
$ws.Range("J141").Value = 99782.13
$ws.Range("L141").Value = 99782.13
$ws.Range("N141").Value = -110142.13

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 50211.11
$ws.Range("J37").Value = 50211.11
$ws.Range("L37").Value = 150633.33
$ws.Range("N37").Value = -150857.33
$ws.Range("H70").Value = 1900
$ws.Range("I70").Value = 1900
$ws.Range("K70").Value = 5700
$ws.Range("M70").Value = -5385
$ws.Range("H73").Value = 1900
$ws.Range("I73").Value = 1900
$ws.Range("K73").Value = 5700
$ws.Range("M73").Value = -4608
$ws.Range("H80").Value = 3749
$ws.Range("I80").Value = 3332.3333
$ws.Range("J80").Value = 4999
$ws.Range("K80").Value = 9996.999899999999
$ws.Range("L80").Value = 14997
$ws.Range("M80").Value = -9060.999899999999
$ws.Range("N80").Value = -16869
$ws.Range("H81").Value = 4959
$ws.Range("J81").Value = 5688.6523
$ws.Range("L81").Value = 17065.9569
$ws.Range("N81").Value = -19311.9569
$ws.Range("H83").Value = 3749
$ws.Range("I83").Value = 3332.3333
$ws.Range("J83").Value = 4999
$ws.Range("K83").Value = 29990.9997
$ws.Range("L83").Value = 44991
$ws.Range("M83").Value = -25310.9997
$ws.Range("N83").Value = -54351
$ws.Range("H84").Value = 4959
$ws.Range("J84").Value = 5688.6523
$ws.Range("L84").Value = 51197.8707
$ws.Range("N84").Value = -62429.8707
$ws.Range("H98").Value = 962.0909
$ws.Range("I98").Value = 713.5
$ws.Range("J98").Value = 1260.4
$ws.Range("K98").Value = 2140.5
$ws.Range("L98").Value = 3781.2
$ws.Range("M98").Value = -642.5
$ws.Range("N98").Value = -6777.200000000001
$ws.Range("H132").Value = 2189.3635
$ws.Range("I132").Value = 1879.2
$ws.Range("J132").Value = 2447.8333
$ws.Range("K132").Value = 16912.8
$ws.Range("L132").Value = 22030.4997
$ws.Range("M132").Value = -14382.8
$ws.Range("N132").Value = -27090.4997

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5230059
$ws.Range("I102").Value = 10101893
$ws.Range("K102").Value = 10101893
$ws.Range("M102").Value = -10100271
$ws.Range("H122").Value = 498482.88
$ws.Range("I122").Value = 687753.25
$ws.Range("J122").Value = 6380
$ws.Range("K122").Value = 2063259.75
$ws.Range("L122").Value = 19140
$ws.Range("M122").Value = -2060809.75
$ws.Range("N122").Value = -24040
$ws.Range("H123").Value = 99107
$ws.Range("J123").Value = 99107
$ws.Range("L123").Value = 99107
$ws.Range("N123").Value = -104007
$ws.Range("H127").Value = 85362.86
$ws.Range("J127").Value = 85362.86
$ws.Range("L127").Value = 85362.86
$ws.Range("N127").Value = -95282.86

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 39725.74
$ws.Range("J22").Value = 1236.5714
$ws.Range("L22").Value = 1236.5714
$ws.Range("N22").Value = -1826.5714
$ws.Range("H27").Value = 39725.74
$ws.Range("J27").Value = 1236.5714
$ws.Range("L27").Value = 1236.5714
$ws.Range("N27").Value = -1450.5714
$ws.Range("H136").Value = 51073.727
$ws.Range("I136").Value = 158228.23
$ws.Range("K136").Value = 474684.6900000001
$ws.Range("M136").Value = -472134.6900000001
$ws.Range("H137").Value = 60933.332
$ws.Range("J137").Value = 60933.332
$ws.Range("L137").Value = 60933.332
$ws.Range("N137").Value = -71133.33199999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 77437.5
$ws.Range("J46").Value = 77437.5
$ws.Range("L46").Value = 77437.5
$ws.Range("N46").Value = -77899.5
$ws.Range("H110").Value = 40000
$ws.Range("J110").Value = 40000
$ws.Range("L110").Value = 40000
$ws.Range("N110").Value = -48180
$ws.Range("H132").Value = 25925186
$ws.Range("I132").Value = 35720348
$ws.Range("J132").Value = 992047.0600000001
$ws.Range("K132").Value = 107161044
$ws.Range("L132").Value = 2976141.18
$ws.Range("M132").Value = -107158514
$ws.Range("N132").Value = -2981201.18
$ws.Range("H134").Value = 77437.5
$ws.Range("J134").Value = 77437.5
$ws.Range("L134").Value = 232312.5
$ws.Range("N134").Value = -237382.5
$ws.Range("H136").Value = 952.35596
$ws.Range("I136").Value = 863.7308
$ws.Range("K136").Value = 2591.1924
$ws.Range("M136").Value = -41.19239999999991
